# "ran code and data" -- relabel the anaerobic/aerobic/O2-limit/fodder-yeast
# conditions on Sheet1 to the "lacto" variants, unify column I (now "H20")
# with column J's value/format, and move the active selection to G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Relabel columns D, E, F, G (rows 2-5 = non-blank set, rows 6-9 = _Blk set) ---
# Written in this order so new shared strings land in the same sequence as the
# target workbook: D-pair, E-base, F-base, E-blk, F-blk, G-pair.
$ws.Range("D2:D5").Value = " Anaerobic lacto"
$ws.Range("D6:D9").Value = " Anaerobic lacto_Blk"

$ws.Range("E2:E5").Value = " Aerobic lacto"
$ws.Range("F2:F5").Value = "O2 limit lacto"

$ws.Range("E6:E9").Value = " Aerobic lacto_Blk"
$ws.Range("F6:F9").Value = "O2 limit lacto_Blk"

$ws.Range("G2:G5").Value = "Fodder yeast"
$ws.Range("G6:G9").Value = "Fodder yeast_Blk"

# --- Column I now mirrors column J (value "H20" + J's cell format) ---
$ws.Range("J2:J9").Copy($ws.Range("I2:I9"))

# --- Move the active selection from G13 to G2 ---
[void]$ws.Range("G2").Select()
